$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.919.88"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "3.408.17"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "409.02"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.77"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +6.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +5.72%  "
$ws.Range("E10").Value = "  +2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.89"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000221"
$ws.Range("E12").Value = "  +33.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.32"
$ws.Range("E13").Value = "  +9.99%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "3.948.90"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.23"
$ws.Range("E16").Value = "  +7.06%  "
$ws.Range("D17").Value = "3.400.81"
$ws.Range("E17").Value = "  -0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.54"
$ws.Range("E18").Value = "  +8.76%  "
$ws.Range("E19").Value = "  +7.00%  "
$ws.Range("D20").Value = "61.914.08"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "445.11"
$ws.Range("E21").Value = "  +41.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "92.08"
$ws.Range("E22").Value = "  +8.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.20"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.18"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.28"
$ws.Range("E25").Value = "  +3.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.39"
$ws.Range("E26").Value = "  +14.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "33.05"
$ws.Range("E27").Value = "  +11.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.79"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.78"
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.99"
$ws.Range("E31").Value = "  +5.13%  "
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.86"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +4.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.78"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  +7.49%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.321"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "143.40"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.29"
$ws.Range("E44").Value = "  +9.84%  "
$ws.Range("E45").Value = "  +15.29%  "
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.151"
$ws.Range("E48").Value = "  +24.81%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.26"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.15"
$ws.Range("E50").Value = "  +7.85%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "3.750.53"
$ws.Range("E51").Value = "  -0.81%  "
